$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column R (18th column) width change: 23.7109375 -> 22.7109375
$ws.Columns.Item(18).ColumnWidth = 22.7109375

# Updated cell values (recalculated ratios/errors)
$ws.Cells.Item(2, 4).Value = [double]"0.007144296897712548"
$ws.Cells.Item(2, 5).Value = [double]"0.007417903189274322"
$ws.Cells.Item(2, 8).Value = [double]"2.083253286508519"
$ws.Cells.Item(2, 9).Value = [double]"0.6590415179803142"
$ws.Cells.Item(2, 10).Value = [double]"0.007567191028275208"
$ws.Cells.Item(2, 11).Value = [double]"0.1481087536733177"
$ws.Cells.Item(2, 12).Value = [double]"5.488204341624448E-05"
$ws.Cells.Item(2, 13).Value = [double]"0.1481087536733174"
$ws.Cells.Item(2, 14).Value = [double]"1.325488815396536"
$ws.Cells.Item(2, 15).Value = [double]"0.3215822963079741"
$ws.Cells.Item(2, 18).Value = [double]"0.2425649162228197"
$ws.Cells.Item(2, 19).Value = [double]"0.3269535019506545"
$ws.Cells.Item(3, 6).Value = [double]"5.185689931824315"
$ws.Cells.Item(3, 7).Value = [double]"0.005759621605348021"
$ws.Cells.Item(3, 14).Value = [double]"0.003275356462361479"
$ws.Cells.Item(3, 15).Value = [double]"2.119924914268915"
$ws.Cells.Item(3, 16).Value = [double]"0.09733120204652208"
$ws.Cells.Item(3, 17).Value = [double]"0.2060866016614308"
$ws.Cells.Item(3, 18).Value = [double]"0.0003203543385384721"
$ws.Cells.Item(3, 19).Value = [double]"2.020137790100577"
$ws.Cells.Item(3, 22).Value = [double]"0.001398585528098584"
$ws.Cells.Item(3, 23).Value = [double]"0.005759715337626024"
$ws.Cells.Item(4, 2).Value = [double]"0.0100170998366197"
$ws.Cells.Item(4, 3).Value = [double]"0.8474172615102259"
$ws.Cells.Item(4, 6).Value = [double]"2.732303424305919"
$ws.Cells.Item(4, 7).Value = [double]"0.01120578579880862"
$ws.Cells.Item(4, 8).Value = [double]"2.073822321355093"
$ws.Cells.Item(4, 9).Value = [double]"0.8201294338742942"
$ws.Cells.Item(4, 10).Value = [double]"0.007564028248501297"
$ws.Cells.Item(4, 11).Value = [double]"0.1745364945928616"
$ws.Cells.Item(4, 12).Value = [double]"5.485910494195209E-05"
$ws.Cells.Item(4, 13).Value = [double]"0.1745364945928622"
$ws.Cells.Item(4, 14).Value = [double]"1.323289409733625"
$ws.Cells.Item(4, 15).Value = [double]"0.2736660921300703"
$ws.Cells.Item(4, 16).Value = [double]"0.1830152333401183"
$ws.Cells.Item(4, 17).Value = [double]"0.3958132555503516"
$ws.Cells.Item(4, 18).Value = [double]"0.2414261194212968"
$ws.Cells.Item(4, 19).Value = [double]"0.3636864051450315"
$ws.Cells.Item(4, 20).Value = [double]"2.658216127534347E-05"
$ws.Cells.Item(4, 21).Value = [double]"0.8486306112154516"
$ws.Cells.Item(4, 22).Value = [double]"0.002654401713710983"
$ws.Cells.Item(4, 23).Value = [double]"0.0112083603976074"
$ws.Cells.Item(4, 24).Value = [double]"0.003665174978805632"
$ws.Cells.Item(4, 25).Value = [double]"0.8486306112154514"
$ws.Cells.Item(5, 14).Value = [double]"0.002818224514931116"
$ws.Cells.Item(5, 15).Value = [double]"1.876895080002627"
$ws.Cells.Item(5, 18).Value = [double]"0.000512323621250908"
$ws.Cells.Item(5, 19).Value = [double]"1.970959586851783"
$ws.Cells.Item(6, 14).Value = [double]"1.32219823094449"
$ws.Cells.Item(6, 15).Value = [double]"0.2209739728222805"
$ws.Cells.Item(6, 18).Value = [double]"0.2396666663667818"
$ws.Cells.Item(6, 19).Value = [double]"0.3656704898636055"
$ws.Cells.Item(7, 2).Value = [double]"0.009976845361061415"
$ws.Cells.Item(7, 3).Value = [double]"0.3297874133506366"
$ws.Cells.Item(7, 14).Value = [double]"0.003639453782854304"
$ws.Cells.Item(7, 15).Value = [double]"1.822070429031099"
$ws.Cells.Item(7, 18).Value = [double]"0.0006632840911626782"
$ws.Cells.Item(7, 19).Value = [double]"1.940759177824897"
$ws.Cells.Item(7, 20).Value = [double]"1.327957135829529E-05"
$ws.Cells.Item(7, 21).Value = [double]"0.3275173631250463"
$ws.Cells.Item(7, 24).Value = [double]"0.001831000578453112"
$ws.Cells.Item(7, 25).Value = [double]"0.3275173631250454"
$ws.Cells.Item(8, 14).Value = [double]"1.319622257679337"
$ws.Cells.Item(8, 15).Value = [double]"0.2325182917524557"
$ws.Cells.Item(8, 18).Value = [double]"0.2382819122323587"
$ws.Cells.Item(8, 19).Value = [double]"0.3670471312898169"
$ws.Cells.Item(9, 2).Value = [double]"0.00998226392451965"
$ws.Cells.Item(9, 3).Value = [double]"0.3367852707192005"
$ws.Cells.Item(9, 8).Value = [double]"4.698652581306738"
$ws.Cells.Item(9, 9).Value = [double]"0.3404243708834113"
$ws.Cells.Item(9, 10).Value = [double]"0.008672415110885429"
$ws.Cells.Item(9, 11).Value = [double]"0.03156690073271559"
$ws.Cells.Item(9, 12).Value = [double]"6.28978257402066E-05"
$ws.Cells.Item(9, 13).Value = [double]"0.03156690073271488"
$ws.Cells.Item(9, 14).Value = [double]"0.003179041188988946"
$ws.Cells.Item(9, 15).Value = [double]"2.1377126504212"
$ws.Cells.Item(9, 18).Value = [double]"0.0003620849407583208"
$ws.Cells.Item(9, 19).Value = [double]"2.044010697369161"
$ws.Cells.Item(9, 20).Value = [double]"1.338794526236375E-05"
$ws.Cells.Item(9, 21).Value = [double]"0.3362066583493737"
$ws.Cells.Item(9, 24).Value = [double]"0.001845943280719976"
$ws.Cells.Item(9, 25).Value = [double]"0.3362066583493737"
$ws.Cells.Item(10, 10).Value = [double]"0.007567712224173581"
$ws.Cells.Item(10, 11).Value = [double]"0.1816623263845749"
$ws.Cells.Item(10, 12).Value = [double]"5.488582345771776E-05"
$ws.Cells.Item(10, 13).Value = [double]"0.181662326384575"
$ws.Cells.Item(10, 14).Value = [double]"1.317214722714466"
$ws.Cells.Item(10, 15).Value = [double]"0.2582171694122649"
$ws.Cells.Item(10, 18).Value = [double]"0.2367852110607417"
$ws.Cells.Item(10, 19).Value = [double]"0.4260308980406651"
$ws.Cells.Item(11, 14).Value = [double]"0.003396968616918185"
$ws.Cells.Item(11, 15).Value = [double]"2.26423125468954"
$ws.Cells.Item(11, 16).Value = [double]"0.2061483249746852"
$ws.Cells.Item(11, 17).Value = [double]"0.3523134073429333"
$ws.Cells.Item(11, 18).Value = [double]"0.0007049510064866325"
$ws.Cells.Item(11, 19).Value = [double]"2.297564997300947"
$ws.Cells.Item(12, 4).Value = [double]"0.007145697660176736"
$ws.Cells.Item(12, 5).Value = [double]"0.006052072924899333"
$ws.Cells.Item(12, 14).Value = [double]"1.317266026380881"
$ws.Cells.Item(12, 15).Value = [double]"0.245071663026405"
$ws.Cells.Item(12, 18).Value = [double]"0.2384620925124516"
$ws.Cells.Item(12, 19).Value = [double]"0.4561006561394234"
